# This script applies the "Updated symbol list" crypto-price refresh to the
# active worksheet (Sheet1 of cryptos.xlsx). The source site re-scraped the
# coin table: most rows keep their coin/link/volume-rank text but get a
# refreshed Price (column D), the hour stamp in column G moves from 17 to 18
# for every data row, and a new top entry ("One") was inserted at row 10,
# pushing WazirX..BTSEToken each down by one row (their Price values were
# also refreshed in the process).
#
# All target values here are text (not numeric) to match the sheet's
# existing inline-string cell typing, so NumberFormat is forced to "@"
# (Text) before each write and restored to "General" afterwards - this
# keeps Excel from "helpfully" re-typing e.g. "242.34" or "18" as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{ Cell = "D2"; Value = "242.34" },
    @{ Cell = "G2"; Value = "18" },
    @{ Cell = "D3"; Value = "23.01" },
    @{ Cell = "G3"; Value = "18" },
    @{ Cell = "D4"; Value = "5.416" },
    @{ Cell = "G4"; Value = "18" },
    @{ Cell = "D5"; Value = "0.05900" },
    @{ Cell = "G5"; Value = "18" },
    @{ Cell = "D6"; Value = "3.440" },
    @{ Cell = "G6"; Value = "18" },
    @{ Cell = "D7"; Value = "6.538" },
    @{ Cell = "G7"; Value = "18" },
    @{ Cell = "D8"; Value = "0.8103" },
    @{ Cell = "G8"; Value = "18" },
    @{ Cell = "D9"; Value = "0.9496" },
    @{ Cell = "G9"; Value = "18" },
    @{ Cell = "B10"; Value = "One" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" },
    @{ Cell = "D10"; Value = "0.01127" },
    @{ Cell = "E10"; Value = "9OneONEBestin24h" },
    @{ Cell = "G10"; Value = "18" },
    @{ Cell = "B11"; Value = "WazirX" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Cell = "D11"; Value = "0.1424" },
    @{ Cell = "E11"; Value = "10WazirXWRX" },
    @{ Cell = "G11"; Value = "18" },
    @{ Cell = "B12"; Value = "MandalaExchangeToken" },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Cell = "D12"; Value = "0.07432" },
    @{ Cell = "E12"; Value = "11MandalaExchangeTokenMDX" },
    @{ Cell = "G12"; Value = "18" },
    @{ Cell = "B13"; Value = "LiechtensteinCryptoassetsExchange" },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" },
    @{ Cell = "D13"; Value = "0.03283" },
    @{ Cell = "E13"; Value = "12LiechtensteinCryptoassetsExchangeLCX" },
    @{ Cell = "G13"; Value = "18" },
    @{ Cell = "B14"; Value = "BitrueCoin" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Cell = "D14"; Value = "0.03053" },
    @{ Cell = "E14"; Value = "13BitrueCoinBTR" },
    @{ Cell = "G14"; Value = "18" },
    @{ Cell = "B15"; Value = "BitMartToken" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Cell = "D15"; Value = "0.09334" },
    @{ Cell = "E15"; Value = "14BitMartTokenBMX" },
    @{ Cell = "G15"; Value = "18" },
    @{ Cell = "B16"; Value = "MCDex" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb" },
    @{ Cell = "D16"; Value = "3.861" },
    @{ Cell = "E16"; Value = "15MCDexMCB" },
    @{ Cell = "G16"; Value = "18" },
    @{ Cell = "B17"; Value = "BitForexToken" },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Cell = "D17"; Value = "0.001566" },
    @{ Cell = "E17"; Value = "16BitForexTokenBF" },
    @{ Cell = "G17"; Value = "18" },
    @{ Cell = "B18"; Value = "CoinExToken" },
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet" },
    @{ Cell = "D18"; Value = "0.04670" },
    @{ Cell = "E18"; Value = "17CoinExTokenCET" },
    @{ Cell = "G18"; Value = "18" },
    @{ Cell = "B19"; Value = "TigerCash" },
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" },
    @{ Cell = "D19"; Value = "0.005876" },
    @{ Cell = "E19"; Value = "18TigerCashTCH" },
    @{ Cell = "G19"; Value = "18" },
    @{ Cell = "B20"; Value = "BitKan" },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan" },
    @{ Cell = "D20"; Value = "0.001260" },
    @{ Cell = "E20"; Value = "19BitKanKAN" },
    @{ Cell = "G20"; Value = "18" },
    @{ Cell = "B21"; Value = "HotbitToken" },
    @{ Cell = "C21"; Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb" },
    @{ Cell = "D21"; Value = "0.004887" },
    @{ Cell = "E21"; Value = "20HotbitTokenHTB" },
    @{ Cell = "G21"; Value = "18" },
    @{ Cell = "B22"; Value = "NitroEx" },
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx" },
    @{ Cell = "D22"; Value = "0.00006803" },
    @{ Cell = "E22"; Value = "21NitroExNTX" },
    @{ Cell = "G22"; Value = "18" },
    @{ Cell = "B23"; Value = "LEO" },
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" },
    @{ Cell = "D23"; Value = "3.572" },
    @{ Cell = "E23"; Value = "22LEOLEO" },
    @{ Cell = "G23"; Value = "18" },
    @{ Cell = "B24"; Value = "BTSEToken" },
    @{ Cell = "C24"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" },
    @{ Cell = "D24"; Value = "2.128" },
    @{ Cell = "E24"; Value = "23BTSETokenBTSE" },
    @{ Cell = "G24"; Value = "18" },
    @{ Cell = "G25"; Value = "18" },
    @{ Cell = "D26"; Value = "0.1329" },
    @{ Cell = "G26"; Value = "18" },
    @{ Cell = "D27"; Value = "0.0002285" },
    @{ Cell = "G27"; Value = "18" },
    @{ Cell = "G28"; Value = "18" },
    @{ Cell = "G29"; Value = "18" },
    @{ Cell = "G30"; Value = "18" },
    @{ Cell = "G31"; Value = "18" },
    @{ Cell = "G32"; Value = "18" },
    @{ Cell = "G33"; Value = "18" },
    @{ Cell = "G34"; Value = "18" },
    @{ Cell = "G35"; Value = "18" },
    @{ Cell = "G36"; Value = "18" },
    @{ Cell = "G37"; Value = "18" },
    @{ Cell = "G38"; Value = "18" },
    @{ Cell = "G39"; Value = "18" },
    @{ Cell = "D40"; Value = "0.03941" },
    @{ Cell = "G40"; Value = "18" },
    @{ Cell = "D41"; Value = "0.006182" },
    @{ Cell = "G41"; Value = "18" },
    @{ Cell = "G42"; Value = "18" },
    @{ Cell = "D43"; Value = "0.003001" },
    @{ Cell = "G43"; Value = "18" },
    @{ Cell = "D44"; Value = "0.009028" },
    @{ Cell = "G44"; Value = "18" },
    @{ Cell = "D45"; Value = "0.00005212" },
    @{ Cell = "G45"; Value = "18" },
    @{ Cell = "G46"; Value = "18" },
    @{ Cell = "D47"; Value = "0.7003" },
    @{ Cell = "G47"; Value = "18" },
    @{ Cell = "G48"; Value = "18" },
    @{ Cell = "D49"; Value = "0.00002101" },
    @{ Cell = "G49"; Value = "18" },
    @{ Cell = "D50"; Value = "0.0002001" },
    @{ Cell = "G50"; Value = "18" },
    @{ Cell = "G51"; Value = "18" }
)

foreach ($edit in $edits) {
    $cell = $ws.Range($edit.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $edit.Value
    $cell.NumberFormat = "General"
}

Write-Host ("Applied {0} cell updates" -f $edits.Count)
